$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets contain the same rows; update the
# "想去人数" (want-to-go count) values for the two affected events.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 137
    $ws.Range("F4").Value = 91
}
